# Updated with the alternate method of wearout
#
# 1) Re-apply the "center" alignment style to the Sheet1 header cells
#    (E1, and F1:H1 keep their existing centering, F2:H2 gain it too).
# 2) Freeze panes on Sheet1 at column H / row 2 (xSplit=8, ySplit=2).
# 3) Append 8 new simulation rows (161-168) with an alternate "wearout"
#    parameter set (Lambda = 0.5, K = 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header styling -------------------------------------------------
# xlCenter
$xlCenter = -4108

$ws.Range("E1").HorizontalAlignment = $xlCenter
$ws.Range("F1:H1").HorizontalAlignment = $xlCenter
$ws.Range("F2:H2").HorizontalAlignment = $xlCenter

# --- Freeze panes -----------------------------------------------------
# Select the cell just below/right of the freeze boundary (row 2, col H)
# then freeze - this locks rows 1:2 and columns A:H.
$ws.Range("I3").Select()
$excel.ActiveWindow.FreezePanes = $true

# Restore the pane selections to match the target view state.
$ws.Range("I1").Select()
$ws.Range("A3").Select()
$ws.Range("I159").Select()

# --- New data rows (161-168) ------------------------------------------
function Set-SimRow {
    param(
        [int]$Row,
        [double]$A,
        [double]$B,
        [double]$C,
        [string]$D,
        [double]$E,
        [double]$F,
        [double]$G,
        [double]$H
    )

    $ws.Range("A$Row").Value = $A
    $ws.Range("B$Row").Value = $B
    $ws.Range("C$Row").Value = $C
    $ws.Range("D$Row").Value = [double]$D
    $ws.Range("D$Row").NumberFormat = $ws.Range("D3").NumberFormat
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H
}

Set-SimRow -Row 161 -A 0.5 -B 0  -C 4  -D "5.7875370370370305E-4" -E 4 -F 17.912369254273202  -G 9.7713736876360393  -H 7.3548237699282701
Set-SimRow -Row 162 -A 0.5 -B 4  -C 4  -D "5.7875370370370305E-4" -E 4 -F 246.52939291626799   -G 130.748737844918    -H 87.520538552008304
Set-SimRow -Row 163 -A 0.5 -B 8  -C 4  -D "5.7875370370370305E-4" -E 4 -F 612.95055273015498   -G 299.00767520110099  -H 206.26904261756599
Set-SimRow -Row 164 -A 0.5 -B 16 -C 4  -D "5.7875370370370305E-4" -E 4 -F 1230.7247064841599   -G 616.12638409976103  -H 415.25637562513299
Set-SimRow -Row 165 -A 0.5 -B 0  -C 16 -D "5.7875370370370305E-4" -E 4 -F 6.1228669329616299   -G 2.5146342016435401  -H 1.68332685048621
Set-SimRow -Row 166 -A 0.5 -B 16 -C 16 -D "5.7875370370370305E-4" -E 4 -F 338.02622985043598   -G 187.22995565284501  -H 119.508993666872
Set-SimRow -Row 167 -A 0.5 -B 32 -C 16 -D "5.7875370370370305E-4" -E 4 -F 1677.2908399560799   -G 847.08392392084204  -H 542.53143373917305
Set-SimRow -Row 168 -A 0.5 -B 64 -C 16 -D "5.7875370370370305E-4" -E 4 -F 4331.0764879380404   -G 2178.46998154115    -H 1438.3704666875999

# Leave the bottom-right pane's selection on the last "old" row, matching
# the target view state.
$ws.Range("I159").Select()
